# Update cryptos list - GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to be written as text (not auto-converted to a number
    # or date by Excel) while keeping the cell's original "Normal" style.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "71.253.05"
$ws.Range("E2").Value = "  +6.24%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.666.93"
$ws.Range("E3").Value = "  +5.82%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.998"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5 - BNB
Set-TextValue "D5" "594.94"
$ws.Range("E5").Value = "  +1.62%  "

# Row 6 - Solana
Set-TextValue "D6" "195.70"
$ws.Range("E6").Value = "  +4.15%  "

# Row 7 - XRP
Set-TextValue "D7" "0.652"
$ws.Range("E7").Value = "  +3.09%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "3.655.55"
$ws.Range("E8").Value = "  +5.69%  "

# Row 9 - USDC
Set-TextValue "D9" "0.999"
$ws.Range("E9").Value = "  -0.05%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.182"
$ws.Range("E10").Value = "  +6.21%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.679"
$ws.Range("E11").Value = "  +4.89%  "

# Row 12 - Avalanche
Set-TextValue "D12" "58.95"
$ws.Range("E12").Value = "  +4.02%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000296"
$ws.Range("E13").Value = "  +6.32%  "

# Row 14 - Polkadot
Set-TextValue "D14" "10.01"
$ws.Range("E14").Value = "  +6.29%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.235.97"
$ws.Range("E15").Value = "  +5.24%  "

# Row 16 - Chainlink
Set-TextValue "D16" "19.96"
$ws.Range("E16").Value = "  +6.54%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "3.647.92"
$ws.Range("E17").Value = "  +5.03%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "70.991.00"
$ws.Range("E18").Value = "  +5.73%  "

# Row 19 - Uniswap
Set-TextValue "D19" "12.84"
$ws.Range("E19").Value = "  +5.60%  "

# Row 20 - TRON (D unchanged)
$ws.Range("E20").Value = "  +2.64%  "

# Row 21 - Polygon
Set-TextValue "D21" "1.08"
$ws.Range("E21").Value = "  +5.66%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "493.58"
$ws.Range("E22").Value = "  +1.57%  "

# Row 23 - InternetComputer(DFINITY)
Set-TextValue "D23" "19.03"
$ws.Range("E23").Value = "  +12.98%  "

# Row 24 - Toncoin
Set-TextValue "D24" "5.38"
$ws.Range("E24").Value = "  +0.01%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "4.53"
$ws.Range("E25").Value = "  +1.41%  "

# Row 26 - Litecoin
Set-TextValue "D26" "92.14"
$ws.Range("E26").Value = "  +2.75%  "

# Row 27 - ImmutableX
Set-TextValue "D27" "3.19"
$ws.Range("E27").Value = "  +8.52%  "

# Row 28 - RenderToken
Set-TextValue "D28" "11.59"
$ws.Range("E28").Value = "  +5.96%  "

# Row 29 - Filecoin
Set-TextValue "D29" "9.68"
$ws.Range("E29").Value = "  +6.59%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "33.05"
$ws.Range("E30").Value = "  +5.35%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "7.95"
$ws.Range("E31").Value = "  +11.34%  "

# Row 32 - Hedera (D unchanged)
$ws.Range("E32").Value = "  +9.33%  "

# Row 33 - Bittensor
Set-TextValue "D33" "633.39"
$ws.Range("E33").Value = "  +5.28%  "

# Row 34 - Cosmos
Set-TextValue "D34" "12.37"
$ws.Range("E34").Value = "  +5.38%  "

# Row 35 - OKB
Set-TextValue "D35" "65.84"
$ws.Range("E35").Value = "  +2.48%  "

# Row 36 - InjectiveProtocol
Set-TextValue "D36" "40.81"
$ws.Range("E36").Value = "  +11.51%  "

# Row 37 - PEPE
Set-TextValue "D37" "0.0₃0844"
$ws.Range("E37").Value = "  +11.60%  "

# Row 38 - TheGraph
Set-TextValue "D38" "0.416"
$ws.Range("E38").Value = "  +8.52%  "

# Row 39 - Kaspa (D unchanged)
$ws.Range("E39").Value = "  -1.17%  "

# Row 40 - Dai (D unchanged)
$ws.Range("E40").Value = "  +0.05%  "

# Row 41 - Stacks
Set-TextValue "D41" "3.64"
$ws.Range("E41").Value = "  +3.23%  "

# Row 42 - Maker
Set-TextValue "D42" "3.327.44"
$ws.Range("E42").Value = "  +2.89%  "

# Rows 43 and 44 swap places: Fetch.AI now ranks above ThetaToken
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D43" "2.90"
$ws.Range("E43").Value = "  +15.25%  "

$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D44" "3.18"
$ws.Range("E44").Value = "  +9.54%  "

# Row 45 - VeChain
Set-TextValue "D45" "0.0456"
$ws.Range("E45").Value = "  +6.16%  "

# Row 46 - dogwifhat (E unchanged)
Set-TextValue "D46" "2.96"

# Row 47 - ApeXProtocol
Set-TextValue "D47" "3.32"
$ws.Range("E47").Value = "  +1.74%  "

# Row 48 - Stellar (D unchanged)
$ws.Range("E48").Value = "  +3.10%  "

# Row 49 - THORChain
Set-TextValue "D49" "9.27"
$ws.Range("E49").Value = "  +6.57%  "

# Row 50 - LidoDAOToken
Set-TextValue "D50" "3.34"
$ws.Range("E50").Value = "  +2.26%  "

# Row 51 - FirstDigitalUSD
Set-TextValue "D51" "0.999"
$ws.Range("E51").Value = "  -0.14%  "
